$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Every existing "SKIP" result (rows 2-21) is now "PASS"
for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 5).Value = "PASS"
}

# New row 23: TestCase_E22 (Unwatch Article)
$ws.Range("A23").Value = "TestCase_E22"
$ws.Range("C23").Value = "Verify that user is able to unwatch an Article from watchlist page"
$ws.Range("B23").Value = "OPQA-293"

# Descriptions for new rows 24-25 (Unwatch Patent, Unwatch Post)
$ws.Range("C24").Value = "Verify that user is able to unwatch a Patent from watchlist page"
$ws.Range("C25").Value = "Verify that user is able to unwatch a Post from watchlist page"

# TCIDs for new rows 24-25
$ws.Range("A24").Value = "TestCase_E23"
$ws.Range("A25").Value = "TestCase_E24"

# Jira ids for new rows 24-25
$ws.Range("B24").Value = "OPQA-294"
$ws.Range("B25").Value = "OPQA-295"

# Runmode + Results columns for all three new rows
$ws.Range("D23").Value = "Y"
$ws.Range("E23").Value = "PASS"
$ws.Range("D24").Value = "Y"
$ws.Range("E24").Value = "PASS"
$ws.Range("D25").Value = "Y"
$ws.Range("E25").Value = "PASS"

# Match formatting (border, wrap) of the new rows to the existing table rows
$ws.Range("A22:E22").Copy()
$ws.Range("A23:E25").PasteSpecial(-4122)

# Sheet view: scrolled so row 13 is at top, selection on A16
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("A16").Select()

# Workbook window size
$excel.ActiveWindow.Width = 14175
$excel.ActiveWindow.Height = 7860
